$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers I1 / J1 - copy formatting (bold, border, centered) from H1 so they
# reuse the existing header style, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data columns I (I0) and J (IF) for rows 2-10
$data = @{
    2  = @(5, 8)
    3  = @(1, 4)
    4  = @(1, 5)
    5  = @(1, 4)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(6, 8)
    9  = @(6, 9)
    10 = @(5, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}
